$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Archive": the localization-status report is refreshed
# for rows that are now awaiting translation. Every cell that previously read
# "Ready for handoff" is updated to "In Translation" (the language status
# columns on the Overview sheet, and the Status column on each per-language
# sheet). The status columns are then resized to fit the new, shorter text.
# ---------------------------------------------------------------------------

$newStatus = "In Translation"

# --- "Overview" sheet: columns E (zh-cn) and F (de-de), rows 2-3 -----------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- "zh-cn" sheet: "Status" column C, rows 2-3 -----------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# --- "de-de" sheet: "Status" column C, rows 2-3 -----------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Resize the status columns to fit the shorter "In Translation" text ----
# (the shared text is narrower than the old "Ready for handoff", so the
# columns shrink to match).
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
